$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E12").Value = 107
$ws.Range("F12").Value = 108
$ws.Range("E12:F12").Style = "Percent"
